$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update bus voltage magnitude results (vm_pu) for Case_3_29 (380 kV case, Vm_pu=1.02)
# Rows 2-25, columns B:F and I:N (G stays 1, H stays empty).

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.059748109965638
$ws.Range("D2").Value = 1.063861502042142
$ws.Range("E2").Value = 1.055552643545337
$ws.Range("F2").Value = 1.074316816113762
$ws.Range("I2").Value = 1.055531757346504
$ws.Range("J2").Value = 1.064732532764985
$ws.Range("K2").Value = 1.066578996682135
$ws.Range("L2").Value = 1.058292786578198
$ws.Range("M2").Value = 1.077006349704311
$ws.Range("N2").Value = 1.06624457508804
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.060967609431931
$ws.Range("D3").Value = 1.06483651508921
$ws.Range("E3").Value = 1.056606006188037
$ws.Range("F3").Value = 1.075448165908222
$ws.Range("I3").Value = 1.055974796120894
$ws.Range("J3").Value = 1.065604119867755
$ws.Range("K3").Value = 1.067368701774916
$ws.Range("L3").Value = 1.059159001432654
$ws.Range("M3").Value = 1.077954029828008
$ws.Range("N3").Value = 1.067117399944468
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.061756335355835
$ws.Range("D4").Value = 1.065467050076784
$ws.Range("E4").Value = 1.057287394046781
$ws.Range("F4").Value = 1.076180224621914
$ws.Range("I4").Value = 1.056260019171712
$ws.Range("J4").Value = 1.066167180539155
$ws.Range("K4").Value = 1.067878710651292
$ws.Range("L4").Value = 1.059718690861036
$ws.Range("M4").Value = 1.078566651344324
$ws.Range("N4").Value = 1.067681260226605
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.062087828837681
$ws.Range("D5").Value = 1.065732040931018
$ws.Range("E5").Value = 1.057573800862858
$ws.Range("F5").Value = 1.076487983252342
$ws.Range("I5").Value = 1.056379579894472
$ws.Range("J5").Value = 1.066403673578297
$ws.Range("K5").Value = 1.068092884068946
$ws.Range("L5").Value = 1.059953791662008
$ws.Range("M5").Value = 1.07882405709579
$ws.Range("N5").Value = 1.067918089112989
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.062143483017544
$ws.Range("D6").Value = 1.065776529044349
$ws.Range("E6").Value = 1.057621887010309
$ws.Range("F6").Value = 1.076539657382974
$ws.Range("I6").Value = 1.056399634335864
$ws.Range("J6").Value = 1.066443369072638
$ws.Range("K6").Value = 1.06812883101835
$ws.Range("L6").Value = 1.05999325485858
$ws.Range("M6").Value = 1.078867268457598
$ws.Range("N6").Value = 1.067957840979486
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.061760765126768
$ws.Range("D7").Value = 1.065470591234429
$ws.Range("E7").Value = 1.057291221218287
$ws.Range("F7").Value = 1.076184336901777
$ws.Range("I7").Value = 1.056261618110701
$ws.Range("J7").Value = 1.06617034142443
$ws.Range("K7").Value = 1.067881573366366
$ws.Range("L7").Value = 1.059721833044916
$ws.Range("M7").Value = 1.078570091362582
$ws.Range("N7").Value = 1.067684425600701
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.06016032376732
$ws.Range("D8").Value = 1.064191088168112
$ws.Range("E8").Value = 1.055908676219138
$ws.Range("F8").Value = 1.074699161854333
$ws.Range("I8").Value = 1.055681785592648
$ws.Range("J8").Value = 1.065027279810751
$ws.Range("K8").Value = 1.066846085409447
$ws.Range("L8").Value = 1.058585696257073
$ws.Range("M8").Value = 1.077326745666483
$ws.Range("N8").Value = 1.06653974070842
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.057337199689381
$ws.Range("D9").Value = 1.061933609730102
$ws.Range("E9").Value = 1.053470806309265
$ws.Range("F9").Value = 1.072082037187428
$ws.Range("I9").Value = 1.054648888420612
$ws.Range("J9").Value = 1.06300599886113
$ws.Range("K9").Value = 1.06501384650766
$ws.Range("L9").Value = 1.056577425894221
$ws.Range("M9").Value = 1.075131235306018
$ws.Range("N9").Value = 1.064515589308005
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.055453004611754
$ws.Range("D10").Value = 1.060426654943902
$ws.Range("E10").Value = 1.051844372554451
$ws.Range("F10").Value = 1.070337170529974
$ws.Range("I10").Value = 1.053952741289574
$ws.Range("J10").Value = 1.061653649961865
$ws.Range("K10").Value = 1.063787191521483
$ws.Range("L10").Value = 1.055234300077809
$ws.Range("M10").Value = 1.073664411357265
$ws.Range("N10").Value = 1.063161319918187
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.054636595727571
$ws.Range("D11").Value = 1.059773644327681
$ws.Range("E11").Value = 1.05113980896919
$ws.Range("F11").Value = 1.069581577886335
$ws.Range("I11").Value = 1.053649500489079
$ws.Range("J11").Value = 1.061066903995405
$ws.Range("K11").Value = 1.063254796091382
$ws.Range("L11").Value = 1.054651679445574
$ws.Range("M11").Value = 1.073028496066345
$ws.Range("N11").Value = 1.062573740705154
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.054333261259951
$ws.Range("D12").Value = 1.059531012170665
$ws.Range("E12").Value = 1.050878054885325
$ws.Range("F12").Value = 1.069300907356367
$ws.Range("I12").Value = 1.053536591239572
$ws.Range("J12").Value = 1.06084878293384
$ws.Range("K12").Value = 1.063056852111848
$ws.Range("L12").Value = 1.054435110616677
$ws.Range("M12").Value = 1.072792171431763
$ws.Range("N12").Value = 1.062355309886662
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.054398331371588
$ws.Range("D13").Value = 1.059583061006214
$ws.Range("E13").Value = 1.050934204214649
$ws.Range("F13").Value = 1.069361112607479
$ws.Range("I13").Value = 1.053560822986838
$ws.Range("J13").Value = 1.060895578701038
$ws.Range("K13").Value = 1.063099320325893
$ws.Range("L13").Value = 1.05448157252231
$ws.Range("M13").Value = 1.072842869205137
$ws.Range("N13").Value = 1.062402172109217
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.054611523732333
$ws.Range("D14").Value = 1.059753589814939
$ws.Range("E14").Value = 1.051118173279688
$ws.Range("F14").Value = 1.069558377782042
$ws.Range("I14").Value = 1.053640172932185
$ws.Range("J14").Value = 1.061048877670954
$ws.Range("K14").Value = 1.063238437836342
$ws.Range("L14").Value = 1.054633781022531
$ws.Range("M14").Value = 1.073008963801978
$ws.Range("N14").Value = 1.062555688781254
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.054742867388315
$ws.Range("D15").Value = 1.059858648269027
$ws.Range("E15").Value = 1.051231516308762
$ws.Range("F15").Value = 1.069679918020266
$ws.Range("I15").Value = 1.053689026964389
$ws.Range("J15").Value = 1.061143306656148
$ws.Range("K15").Value = 1.063324127689479
$ws.Range("L15").Value = 1.05472754077378
$ws.Range("M15").Value = 1.073111284563657
$ws.Range("N15").Value = 1.062650251866438
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.05550717537081
$ws.Range("D16").Value = 1.060469982677601
$ws.Range("E16").Value = 1.051891125507931
$ws.Range("F16").Value = 1.070387315446506
$ws.Range("I16").Value = 1.053972828278008
$ws.Range("J16").Value = 1.061692565563476
$ws.Range("K16").Value = 1.063822498482594
$ws.Range("L16").Value = 1.055272944678832
$ws.Range("M16").Value = 1.073706598586032
$ws.Range("N16").Value = 1.063200290784416
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.055986459079758
$ws.Range("D17").Value = 1.06085332460312
$ws.Range("E17").Value = 1.052304797552179
$ws.Range("F17").Value = 1.070831031551544
$ws.Range("I17").Value = 1.054150365516448
$ws.Range("J17").Value = 1.06203678672431
$ws.Range("K17").Value = 1.064134778771059
$ws.Range("L17").Value = 1.055614783034823
$ws.Range("M17").Value = 1.074079815859714
$ws.Range("N17").Value = 1.06354500077878
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.056265965243832
$ws.Range("D18").Value = 1.061076874355474
$ws.Range("E18").Value = 1.05254605594634
$ws.Range("F18").Value = 1.071089838566604
$ws.Range("I18").Value = 1.054253745944798
$ws.Range("J18").Value = 1.062237452335332
$ws.Range("K18").Value = 1.064316806484873
$ws.Range("L18").Value = 1.05581407149325
$ws.Range("M18").Value = 1.074297432913522
$ws.Range("N18").Value = 1.063745951357987
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.056361260899172
$ws.Range("D19").Value = 1.061153091096668
$ws.Range("E19").Value = 1.052628313868544
$ws.Range("F19").Value = 1.071178084259039
$ws.Range("I19").Value = 1.054288966509467
$ws.Range("J19").Value = 1.06230585504684
$ws.Range("K19").Value = 1.064378852904659
$ws.Range("L19").Value = 1.05588200678145
$ws.Range("M19").Value = 1.074371622155104
$ws.Range("N19").Value = 1.063814451209191
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.055935041864429
$ws.Range("D20").Value = 1.060812200526246
$ws.Range("E20").Value = 1.05226041752226
$ws.Range("F20").Value = 1.07078342555165
$ws.Range("I20").Value = 1.054131335460507
$ws.Range("J20").Value = 1.061999866728875
$ws.Range("K20").Value = 1.064101286484863
$ws.Range("L20").Value = 1.055578117377606
$ws.Range("M20").Value = 1.074039780872492
$ws.Range("N20").Value = 1.063508028352718
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.054548746232609
$ws.Range("D21").Value = 1.059703375411313
$ws.Range("E21").Value = 1.051064000308214
$ws.Range("F21").Value = 1.069500288388366
$ws.Range("I21").Value = 1.053616813870824
$ws.Range("J21").Value = 1.061003739858391
$ws.Range("K21").Value = 1.063197476415245
$ws.Range("L21").Value = 1.05458896377584
$ws.Range("M21").Value = 1.072960056335947
$ws.Range("N21").Value = 1.06251048686782
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.053676640617779
$ws.Range("D22").Value = 1.059005778765292
$ws.Range("E22").Value = 1.050311487948383
$ws.Range("F22").Value = 1.068693471707136
$ws.Range("I22").Value = 1.053291738855502
$ws.Range("J22").Value = 1.06037640763712
$ws.Range("K22").Value = 1.062628123178868
$ws.Range("L22").Value = 1.053966131118011
$ws.Range("M22").Value = 1.072280510752014
$ws.Range("N22").Value = 1.061882263762839
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.054139007335986
$ws.Range("D23").Value = 1.059375629592312
$ws.Range("E23").Value = 1.050710435838605
$ws.Range("F23").Value = 1.069121186502729
$ws.Range("I23").Value = 1.053464216888673
$ws.Range("J23").Value = 1.060709066305567
$ws.Range("K23").Value = 1.062930052095641
$ws.Range("L23").Value = 1.054296393524692
$ws.Range("M23").Value = 1.072640815697431
$ws.Range("N23").Value = 1.062215394844752
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.055958275251212
$ws.Range("D24").Value = 1.060830782872062
$ws.Range("E24").Value = 1.052280471038725
$ws.Range("F24").Value = 1.070804936669621
$ws.Range("I24").Value = 1.054139934861423
$ws.Range("J24").Value = 1.062016549634083
$ws.Range("K24").Value = 1.064116420578459
$ws.Range("L24").Value = 1.055594685319101
$ws.Range("M24").Value = 1.074057871189475
$ws.Range("N24").Value = 1.063524734949564
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.058067407389227
$ws.Range("D25").Value = 1.06251756423756
$ws.Range("E25").Value = 1.054101257594447
$ws.Range("F25").Value = 1.072758641447195
$ws.Range("I25").Value = 1.054917244223051
$ws.Range("J25").Value = 1.063529393481127
$ws.Range("K25").Value = 1.065488429080743
$ws.Range("L25").Value = 1.057097360579136
$ws.Range("M25").Value = 1.075699377236589
$ws.Range("N25").Value = 1.065039727208397
